$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.205.68'
$ws.Range('E2').Value = '  -0.32%  '
$ws.Range('D3').Value = '2.641.06'
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'596.94"
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').Value = "'158.95"
$ws.Range('E6').Value = '  +2.66%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = "'0.543"
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('D9').Value = "'0.142"
$ws.Range('E9').Value = '  +2.51%  '
$ws.Range('D10').Value = "'0.156"
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('D11').Value = "'5.26"
$ws.Range('E11').Value = '  +0.17%  '
$ws.Range('D12').Value = "'0.350"
$ws.Range('E12').Value = '  -0.60%  '
$ws.Range('D13').Value = "'27.92"
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').Value = "'0.0000188"
$ws.Range('E14').Value = '  +0.13%  '
$ws.Range('D15').Value = '3.119.43'
$ws.Range('E15').Value = '  -0.05%  '
$ws.Range('D16').Value = '68.093.82'
$ws.Range('E16').Value = '  -0.32%  '
$ws.Range('D17').Value = '2.620.87'
$ws.Range('E17').Value = '  -0.88%  '
$ws.Range('D18').Value = "'11.39"
$ws.Range('E18').Value = '  -1.53%  '
$ws.Range('D19').Value = "'359.94"
$ws.Range('E19').Value = '  -1.69%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = "'7.32"
$ws.Range('E20').Value = '  -2.26%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').Value = "'4.39"
$ws.Range('E21').Value = '  +0.70%  '
$ws.Range('D22').Value = "'4.77"
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D23').Value = "'2.05"
$ws.Range('E23').Value = '  -2.87%  '
$ws.Range('D24').Value = "'74.93"
$ws.Range('E24').Value = '  +1.87%  '
$ws.Range('E25').Value = '  -0.08%  '
$ws.Range('E26').Value = '  -2.44%  '
$ws.Range('B27').Value = 'WrappedeETH'
$ws.Range('C27').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D27').Value = '2.777.02'
$ws.Range('E27').Value = '  +0.22%  '
$ws.Range('B28').Value = 'PEPE'
$ws.Range('C28').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D28').Value = "'0.0000104"
$ws.Range('E28').Value = '  -0.36%  '
$ws.Range('B29').Value = 'Bittensor'
$ws.Range('C29').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D29').Value = "'559.35"
$ws.Range('E29').Value = '  -2.92%  '
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = "'0.894"
$ws.Range('E30').Value = '  -10.55%  '
$ws.Range('D31').Value = "'7.99"
$ws.Range('E31').Value = '  -0.93%  '
$ws.Range('D32').Value = "'1.39"
$ws.Range('E32').Value = '  -2.70%  '
$ws.Range('D33').Value = "'1.86"
$ws.Range('E33').Value = '  -0.85%  '
$ws.Range('B34').Value = 'FirstDigitalUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D34').Value = "'0.999"
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('B35').Value = 'Kaspa'
$ws.Range('C35').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D35').Value = "'0.128"
$ws.Range('E35').Value = '  -2.45%  '
$ws.Range('D36').Value = "'1.55"
$ws.Range('E36').Value = '  -2.66%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').Value = "'160.87"
$ws.Range('E37').Value = '  +1.20%  '
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').Value = "'19.73"
$ws.Range('E38').Value = '  +1.42%  '
$ws.Range('D39').Value = "'0.370"
$ws.Range('E39').Value = '  -0.25%  '
$ws.Range('D40').Value = "'1.87"
$ws.Range('E40').Value = '  -2.59%  '
$ws.Range('D41').Value = "'5.33"
$ws.Range('E41').Value = '  -1.96%  '
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').Value = "'17.79"
$ws.Range('E42').Value = '  +0.32%  '
$ws.Range('D43').Value = "'2.64"
$ws.Range('E43').Value = '  -1.58%  '
$ws.Range('B44').Value = 'BabyDogeCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D44').Value = '0.0₆0330'
$ws.Range('E44').Value = '  +2.48%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('D46').Value = "'157.35"
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('D47').Value = "'3.76"
$ws.Range('E47').Value = '  -0.86%  '
$ws.Range('D48').Value = "'22.19"
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('D49').Value = "'1.68"
$ws.Range('E49').Value = '  -2.55%  '
$ws.Range('D50').Value = "'0.0775"
$ws.Range('E50').Value = '  -1.13%  '
$ws.Range('D51').Value = "'0.615"
$ws.Range('E51').Value = '  -0.25%  '
